$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a copy of row 6 at row 7, pushing existing rows 7-53 down to 8-54.
# This reproduces a new weekly data point being inserted into the series.
$ws.Rows.Item(6).Copy()
$ws.Rows.Item(7).Insert()
